# Ran code for averaged intensities on spiral schemes.
#
# The HKL-index table on GossA gains three new sampling schemes
# ("Spiral-90deg-10rot-5space", "Spiral-90deg-15rot-5space",
# "Spiral-90deg-10rot-3space"). The previously-last scheme
# "Gaussian-Quadrature" is re-ordered to sit right after
# "Ring Perpendicular to TD" (immediately before the new Spiral rows),
# and every scheme that used to follow it shifts down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P")

# --- 1. Snapshot the existing per-scheme rows (10-16) before overwriting ---
# Row 10: NoRotation-tilt60deg   Row 13: HexGrid-90degTilt5degRes
# Row 11: Rotation-NoTilt        Row 14: HexGrid-90degTilt22p5degRes
# Row 12: Rotation-60detTilt     Row 15: HexGrid-60degTilt5degRes
# Row 16: Gaussian-Quadrature
$snapshot = @{}
foreach ($r in 10..16) {
    $rowVals = @{}
    $rowVals["B"] = $ws.Cells.Item($r, 2).Value2
    foreach ($c in $cols) {
        $col = $c
        $colIndex = $ws.Range($col + "1").Column
        $rowVals[$col] = $ws.Cells.Item($r, $colIndex).Value2
    }
    $snapshot[$r] = $rowVals
}

function Write-Row($r, $bValue, $data) {
    $ws.Cells.Item($r, 2).Value2 = $bValue
    foreach ($c in $cols) {
        $colIndex = $ws.Range($c + "1").Column
        $ws.Cells.Item($r, $colIndex).Value2 = $data[$c]
    }
}

# --- 2. New row 10: Gaussian-Quadrature (data formerly on row 16) ---
Write-Row 10 "Gaussian-Quadrature" $snapshot[16]

# --- 3. New rows 11-13: brand-new Spiral scheme data ---
$spiral1 = @{ "C"=0.7736688538060419; "D"=1.303916321476776; "E"=0.9307027969858641; "F"=1.059554205819293; "G"=0.7736688538060419; "H"=1.303916321476776; "I"=0.8801929826685201; "J"=1.063820884710824; "K"=0.913226274997322; "L"=1.196573539604364; "M"=0.7736688538060419; "N"=1.11730955923132; "O"=1.016960544521994; "P"=1.015206982508626 }
$spiral2 = @{ "C"=0.7729313415068712; "D"=1.305311981312687; "E"=0.9302468709536855; "F"=1.059863588466121; "G"=0.7729313415068712; "H"=1.305311981312687; "I"=0.8795977688326166; "J"=1.063991938556708; "K"=0.912915563974888; "L"=1.197517652890035; "M"=0.7729313415068712; "N"=1.117779426133187; "O"=1.017088445559841; "P"=1.015297088311702 }
$spiral3 = @{ "C"=0.7735119461771636; "D"=1.304129990706679; "E"=0.9306547013103774; "F"=1.059609727862718; "G"=0.7735119461771636; "H"=1.304129990706679; "I"=0.880057968614301; "J"=1.063856850715124; "K"=0.9131621906631909; "L"=1.196764343698158; "M"=0.7735119461771636; "N"=1.117392346008528; "O"=1.016976591514235; "P"=1.015218464968464 }

Write-Row 11 "Spiral-90deg-10rot-5space" $spiral1
Write-Row 12 "Spiral-90deg-15rot-5space" $spiral2
Write-Row 13 "Spiral-90deg-10rot-3space" $spiral3

# --- 4. New rows 14-19: the old rows 10-15, shifted down by 4 ---
Write-Row 14 "NoRotation-tilt60deg" $snapshot[10]
Write-Row 15 "Rotation-NoTilt" $snapshot[11]
Write-Row 16 "Rotation-60detTilt" $snapshot[12]

# Rows 17-19 are brand new sheet rows; copy row 3's formatting onto column A
# first so the HKL-index cell picks up the same bold/border/centered style
# (s="1") as every other row, then fill in the values.
$ws.Range("A3").Copy()
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Copy()
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Copy()
$ws.Range("A19").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(17, 1).Value2 = 15
$ws.Cells.Item(18, 1).Value2 = 16
$ws.Cells.Item(19, 1).Value2 = 17

Write-Row 17 "HexGrid-90degTilt5degRes" $snapshot[13]
Write-Row 18 "HexGrid-90degTilt22p5degRes" $snapshot[14]
Write-Row 19 "HexGrid-60degTilt5degRes" $snapshot[15]
